# Combining data quality measurement with data model quality.
#
# Physical worksheet positions (tab order) are preserved; only tab names and
# the data inside sheets 2-5 change (sheet 1, "0_SCHEMA_METADATA", is left
# untouched).
#
#   pos2 (was "1_ISSUES")           -> "SCHEME_MEASURES"    (was issues table 8x18 -> becomes 3x6 measures table)
#   pos3 (was "2_SCHEME_MEASURES")  -> "METADATA_ISSUES"     (was 3x6 measures table -> becomes issues table 8x18)
#   pos4 (was "3_MODEL_MEASURES")   -> "METADATA_MEASURES"   (3x3 -> 3x4 measures table)
#   pos5 (was "4_MODEL_METRICS")    -> "METADATA_METRICS"    (3x9 -> 3x8 metrics table)

$wb = $excel.ActiveWorkbook

$wsIssues   = $wb.Worksheets.Item(2)   # currently "1_ISSUES" (8 cols x 18 rows of data)
$wsMeasures = $wb.Worksheets.Item(3)   # currently "2_SCHEME_MEASURES" (3 cols x 6 rows of data)
$wsModelMeasures = $wb.Worksheets.Item(4)  # currently "3_MODEL_MEASURES"
$wsModelMetrics  = $wb.Worksheets.Item(5)  # currently "4_MODEL_METRICS"

# ---------------------------------------------------------------------
# Step 1: while $wsIssues (pos 2) still has its full 8-column header
# formatting, copy the D1:H1 header cell formatting over to $wsMeasures
# (pos 3) so its new, wider header row (rule/desc/owner/table/column/
# constraint_name/length/limit) keeps the same bold/centered/bordered
# style as columns A1:C1 already have there.
# ---------------------------------------------------------------------
$wsIssues.Range("D1:H1").Copy()
$wsMeasures.Range("D1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Step 2: write the full issues table into $wsMeasures (pos 3), which
# becomes "METADATA_ISSUES". This is the same table that used to live in
# $wsIssues, with the rule codes/description renamed:
#   MQMD06 -> MQME10 (unchanged description)
#   MQMD10 -> MQME01 (description shortened to "Columns without comments")
# ---------------------------------------------------------------------
$wsMeasures.Cells.Item(1,1).Value = "rule"
$wsMeasures.Cells.Item(1,2).Value = "desc"
$wsMeasures.Cells.Item(1,3).Value = "owner"
$wsMeasures.Cells.Item(1,4).Value = "table"
$wsMeasures.Cells.Item(1,5).Value = "column"
$wsMeasures.Cells.Item(1,6).Value = "constraint_name"
$wsMeasures.Cells.Item(1,7).Value = "length"
$wsMeasures.Cells.Item(1,8).Value = "limit"

$issuesRows = @(
    @(2,  "MQME10", "Total number of tables with plural names", "SIPAJ", "PROTOCOLO_INTENCOES", ""),
    @(3,  "MQME10", "Total number of tables with plural names", "SIPAJ", "SOL_CRED_RESS_INFORMA_FIS", ""),
    @(4,  "MQME01", "Columns without comments", "SIPAJ", "ANDAMENTO", "STA_PROCESSO"),
    @(5,  "MQME01", "Columns without comments", "SIPAJ", "DOCUMENTO", "DSC_OBSERVACAO_CASSACAO"),
    @(6,  "MQME01", "Columns without comments", "SIPAJ", "INFORMACAO_FISCAL", "VLR_OPR_INTERESTADUAL"),
    @(7,  "MQME01", "Columns without comments", "SIPAJ", "PROC_ANEXO", "TXT_PROC_ANEXO"),
    @(8,  "MQME01", "Columns without comments", "SIPAJ", "PROC_ANEXO_TMP", "TXT_PROC_ANEXO"),
    @(9,  "MQME01", "Columns without comments", "SIPAJ", "PROC_DOCUMENTO", "DSC_OBSERVACAO"),
    @(10, "MQME01", "Columns without comments", "SIPAJ", "PROC_DOCUMENTO_PENDENCIA", "SEQ_PROC_DOCUMENTO"),
    @(11, "MQME01", "Columns without comments", "SIPAJ", "PROC_DOCUMENTO_PENDENCIA", "STA_PENDENCIA"),
    @(12, "MQME01", "Columns without comments", "SIPAJ", "PROC_DOCUMENTO_20190212", "SEQ_PROC_DOCUMENTO"),
    @(13, "MQME01", "Columns without comments", "SIPAJ", "PROC_DOCUMENTO_20190212", "NUM_PROC_DOCUMENTO"),
    @(14, "MQME01", "Columns without comments", "SIPAJ", "PROC_DOCUMENTO_20190212", "DAT_INCLUSAO"),
    @(15, "MQME01", "Columns without comments", "SIPAJ", "PROC_DOCUMENTO_20190212", "DSC_OBSERVACAO"),
    @(16, "MQME01", "Columns without comments", "SIPAJ", "PROTOCOLO_INTENCOES", "SEQ_PROC_DOCUMENTO"),
    @(17, "MQME01", "Columns without comments", "SIPAJ", "RECURSO", "DSC_TEOR_DESPACHO"),
    @(18, "MQME01", "Columns without comments", "SIPAJ", "RECURSO", "DSC_OBSERVACAO")
)

foreach ($row in $issuesRows) {
    $r = $row[0]
    $wsMeasures.Cells.Item($r,1).Value = $row[1]
    $wsMeasures.Cells.Item($r,2).Value = $row[2]
    $wsMeasures.Cells.Item($r,3).Value = $row[3]
    $wsMeasures.Cells.Item($r,4).Value = $row[4]
    $wsMeasures.Cells.Item($r,5).Value = $row[5]
}

# ---------------------------------------------------------------------
# Step 3: shrink $wsIssues (pos 2) down to the small 3-column/6-row
# measures table it becomes ("SCHEME_MEASURES"). Drop the now-unused
# D:H columns and rows 7:18 first, then overwrite the remaining A1:C6
# with the renamed indicator codes (MQMD0x -> MQMS0x).
# ---------------------------------------------------------------------
$wsIssues.Range("D1:H18").EntireColumn.Delete()
$wsIssues.Range("A7:A18").EntireRow.Delete()

$wsIssues.Cells.Item(1,1).Value = "Indicator"
$wsIssues.Cells.Item(1,2).Value = "Description"
$wsIssues.Cells.Item(1,3).Value = "Value"

$schemeRows = @(
    @(2, "MQMS01", "Total number of tables", 73),
    @(3, "MQMS02", "Total number of columns", 362),
    @(4, "MQMS03", "Total number of primary key", 89),
    @(5, "MQMS04", "Total number of foreign key", 107),
    @(6, "MQMS05", "Total number of unique key", 0)
)

foreach ($row in $schemeRows) {
    $r = $row[0]
    $wsIssues.Cells.Item($r,1).Value = $row[1]
    $wsIssues.Cells.Item($r,2).Value = $row[2]
    $wsIssues.Cells.Item($r,3).Value = $row[3]
}

# ---------------------------------------------------------------------
# Step 4: $wsModelMeasures (pos 4), becomes "METADATA_MEASURES". Row 2
# and 3 get new codes/descriptions/values, and a new row 4 is added.
# ---------------------------------------------------------------------
$wsModelMeasures.Cells.Item(2,1).Value = "MQME00"
$wsModelMeasures.Cells.Item(2,2).Value = "Total number of columns"
$wsModelMeasures.Cells.Item(2,3).Value = 362

$wsModelMeasures.Cells.Item(3,1).Value = "MQMEA1"
$wsModelMeasures.Cells.Item(3,2).Value = "Total number of length-required columns"
$wsModelMeasures.Cells.Item(3,3).Value = 69

$wsModelMeasures.Cells.Item(4,1).Value = "MQMEA2"
$wsModelMeasures.Cells.Item(4,2).Value = "Total number of NUMBER columns"
$wsModelMeasures.Cells.Item(4,3).Value = 235

# ---------------------------------------------------------------------
# Step 5: $wsModelMetrics (pos 5), becomes "METADATA_METRICS". Rows 2-8
# get renamed codes/descriptions (some values also change), and row 9
# (IQMD08 / "Table with standard UK prefixes") is removed entirely.
# ---------------------------------------------------------------------
$wsModelMetrics.Range("A9:C9").EntireRow.Delete()

$metricsRows = @(
    @(2, "IQME01", "Columns with comments", "95.86%"),
    @(3, "IQME02", "Columns with data type", "100.00%"),
    @(4, "IQME03", "Length-required columns with data length", "100.00%"),
    @(5, "IQME04", "NUMBER columns with valid scale", "100.00%"),
    @(6, "IQME05", "Columns with valid num_distinct", "100.00%"),
    @(7, "IQME06", "Columns with valid num_nulls", "100.00%"),
    @(8, "IQME07", "Columns with valid density", "100.00%")
)

foreach ($row in $metricsRows) {
    $r = $row[0]
    $wsModelMetrics.Cells.Item($r,1).Value = $row[1]
    $wsModelMetrics.Cells.Item($r,2).Value = $row[2]
    $wsModelMetrics.Cells.Item($r,3).Value = $row[3]
}

# ---------------------------------------------------------------------
# Step 6: rename the tabs (sheetId / position / r:id stay fixed, only
# the visible name changes).
# ---------------------------------------------------------------------
$wsIssues.Name = "SCHEME_MEASURES"
$wsMeasures.Name = "METADATA_ISSUES"
$wsModelMeasures.Name = "METADATA_MEASURES"
$wsModelMetrics.Name = "METADATA_METRICS"
